$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set D15 to "Y", matching the formatting already used by the other cells
# in column D (e.g. D14): centered horizontally and vertically, no border.
$ws.Range("D15").Value = "Y"
$ws.Range("D15").HorizontalAlignment = -4108
$ws.Range("D15").VerticalAlignment = -4108

# Update the selection/active cell to D15
$ws.Range("D15").Select()
